# Apply the edit described by the diff:
# - Row 1: replace header text labels with numeric index values (0-13), keep style.
# - Row 2: becomes the old header labels (shifted down from row 1), with K2/M2/N2 empty.
# - Rows 3-40: column N gets "18-8 Stainless Steel" (moved from old A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1")

# Row 1: numeric values 0..13 across columns A..N
for ($i = 0; $i -lt 14; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Row 2: old row 1 header labels shift down to row 2.
# K2, M2, N2 are left untouched (they were already empty inlineStr cells and stay that way).
$ws.Cells.Item(2, 1).Value = "Lg."
$ws.Cells.Item(2, 2).Value = "Threading"
$ws.Cells.Item(2, 3).Value = "Min.Thread Lg."
$ws.Cells.Item(2, 4).Value = "HeadDia."
$ws.Cells.Item(2, 5).Value = "Head Ht."
$ws.Cells.Item(2, 6).Value = "DriveSize"
$ws.Cells.Item(2, 7).Value = "Finish"
$ws.Cells.Item(2, 8).Value = "TensileStrength, psi"
$ws.Cells.Item(2, 9).Value = "Specifications Met"
$ws.Cells.Item(2, 10).Value = "Pkg.Qty."
$ws.Cells.Item(2, 12).Value = "Pkg."

# Column N for rows 3-40: fill with "18-8 Stainless Steel" (moved from old A2)
for ($r = 3; $r -le 40; $r++) {
    $ws.Cells.Item($r, 14).Value = "18-8 Stainless Steel"
}
